$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: wipe existing cell content (keep formatting) across the working area ---
$ws.Range("A1:AD30").ClearContents()

# --- Step 2: fully remove the now-unused trailing columns (U:AD) for rows 1-2, including their style ---
$ws.Range("U1:AD2").Clear()

# --- Step 3: extend the bold/border/center style used in column A down to the 4 new rows (20-23) ---
$ws.Range("A19").Copy()
$ws.Range("A20:A23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 4: Row 1 numeric column headers (0-18) ---
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 4
$ws.Range("G1").Value = 5
$ws.Range("H1").Value = 6
$ws.Range("I1").Value = 7
$ws.Range("J1").Value = 8
$ws.Range("K1").Value = 9
$ws.Range("L1").Value = 10
$ws.Range("M1").Value = 11
$ws.Range("N1").Value = 12
$ws.Range("O1").Value = 13
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("R1").Value = 16
$ws.Range("S1").Value = 17
$ws.Range("T1").Value = 18

# --- Step 5: Row 2 "HKL" label ---
$ws.Range("B2").Value = "HKL"

# --- Step 6: column B dataset names, in the exact write order needed to rebuild the shared-string
#     table in the same order as the target workbook (Holden entries were spliced into the source
#     list right after the Offset entries, ahead of the HexGrid entries, even though their rows
#     end up physically below the HexGrid rows on the sheet) ---
$ws.Range("B3").Value = "BT8Hex_2.5"
$ws.Range("B4").Value = "BT8Hex_5"
$ws.Range("B5").Value = "BT8Hex_10"
$ws.Range("B6").Value = "BT8Hex_15"
$ws.Range("B7").Value = "Spiral2.5"
$ws.Range("B8").Value = "Spiral5"
$ws.Range("B9").Value = "Spiral7.5"
$ws.Range("B10").Value = "Spiral10"
$ws.Range("B11").Value = "Spiral15"
$ws.Range("B12").Value = "OffsetF45"
$ws.Range("B13").Value = "OffsetA45"
$ws.Range("B14").Value = "OffsetFTD"
$ws.Range("B15").Value = "OffsetATD"
$ws.Range("B20").Value = "Holden2.5"
$ws.Range("B21").Value = "Holden5"
$ws.Range("B22").Value = "Holden10"
$ws.Range("B23").Value = "Holden15"
$ws.Range("B16").Value = "HexGrid-90degTilt2.5degRes"
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("B18").Value = "HexGrid-90degTilt10degRes"
$ws.Range("B19").Value = "HexGrid-90degTilt15degRes"

# --- Step 7: row 2 [h,k,l] headers C2:J2, in their new order ---
$ws.Range("C2").Value = "[2, 1, 1]"
$ws.Range("D2").Value = "[4, 0, 0]"
$ws.Range("E2").Value = "[2, 0, 0]"
$ws.Range("F2").Value = "[2, 2, 0]"
$ws.Range("G2").Value = "[1, 1, 0]"
$ws.Range("H2").Value = "[3, 1, 0]"
$ws.Range("I2").Value = "[2, 2, 2]"
$ws.Range("J2").Value = "[3, 2, 1]"

# --- Step 8: row 2 pair-scheme headers K2:T2 ---
$ws.Range("K2").Value = "1Pair-A"
$ws.Range("L2").Value = "1Pair-B"
$ws.Range("M2").Value = "2Pairs-A"
$ws.Range("N2").Value = "2Pairs-B"
$ws.Range("O2").Value = "3Pairs-A"
$ws.Range("P2").Value = "3Pairs-B"
$ws.Range("Q2").Value = "3Pairs-C"
$ws.Range("R2").Value = "4Pairs"
$ws.Range("S2").Value = "5A4F"
$ws.Range("T2").Value = "MaxUnique"

# --- Step 9: column A row indices (rows 2-23) ---
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16
$ws.Range("A19").Value = 17
$ws.Range("A20").Value = 18
$ws.Range("A21").Value = 19
$ws.Range("A22").Value = 20
$ws.Range("A23").Value = 21

# --- Step 10: data body, every remaining C:T cell for rows 3-23 is 1 ---
$ws.Range("C3:T3").Value = 1
$ws.Range("C4:T4").Value = 1
$ws.Range("C5:T5").Value = 1
$ws.Range("C6:T6").Value = 1
$ws.Range("C7:T7").Value = 1
$ws.Range("C8:T8").Value = 1
$ws.Range("C9:T9").Value = 1
$ws.Range("C10:T10").Value = 1
$ws.Range("C11:T11").Value = 1
$ws.Range("C12:T12").Value = 1
$ws.Range("C13:T13").Value = 1
$ws.Range("C14:T14").Value = 1
$ws.Range("C15:T15").Value = 1
$ws.Range("C16:T16").Value = 1
$ws.Range("C17:T17").Value = 1
$ws.Range("C18:T18").Value = 1
$ws.Range("C19:T19").Value = 1
$ws.Range("C20:T20").Value = 1
$ws.Range("C21:T21").Value = 1
$ws.Range("C22:T22").Value = 1
$ws.Range("C23:T23").Value = 1

Write-Output "Edit applied"
